$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "https://plot.ly/~Dreamshot/9199/import-plotly-plotly-version-/#/"

$ws.Range("A5").Select()
